$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 15.02.2022 09:45"

# 2. D6: change from text "+0.4" to numeric 0.4
$ws.Range("D6").Value = 0.4

# 3. E6: change from text timestamp to a real Excel date/time serial value,
#    matching the style used by the other rows (yyyy-mm-dd hh:mm:ss number format)
$ws.Range("E6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E6").Value = 44607.39802083333
